$wb = $excel.ActiveWorkbook

# --- Step 1: strip header-row styling (bold+border -> Normal) on all sheets ---
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1:N1").Style = "Normal"
}

# --- Step 2: refresh market-price data cells (columns H-N) ---

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 766.3333
$ws.Range("I12").Value = 766.3333
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 766.3333
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -596.3333
$ws.Range("N12").Value = $null
$ws.Range("H40").Value = 5043.2666
$ws.Range("J40").Value = 5062.4165
$ws.Range("L40").Value = 5062.4165
$ws.Range("N40").Value = -5412.4165
$ws.Range("H76").Value = 2756.6924
$ws.Range("J76").Value = 2650
$ws.Range("L76").Value = 2650
$ws.Range("N76").Value = -3280
$ws.Range("H79").Value = 2756.6924
$ws.Range("J79").Value = 2650
$ws.Range("L79").Value = 2650
$ws.Range("N79").Value = -4834
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = $null
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = $null
$ws.Range("H100").Value = 8581.3545
$ws.Range("I100").Value = 4669.909
$ws.Range("K100").Value = 4669.909
$ws.Range("M100").Value = -4128.909
$ws.Range("H106").Value = 101212.73
$ws.Range("I106").Value = 149477.14
$ws.Range("K106").Value = 149477.14
$ws.Range("M106").Value = -148846.14
$ws.Range("H113").Value = 7057.0713
$ws.Range("I113").Value = 5892
$ws.Range("K113").Value = 5892
$ws.Range("M113").Value = -2638

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5783.3076
$ws.Range("J2").Value = 4000
$ws.Range("L2").Value = 4000
$ws.Range("N2").Value = -4226
$ws.Range("H11").Value = 9999
$ws.Range("J11").Value = 9999
$ws.Range("L11").Value = 9999
$ws.Range("N11").Value = -10287
$ws.Range("H45").Value = 3133.375
$ws.Range("I45").Value = 2929.5908
$ws.Range("K45").Value = 2929.5908
$ws.Range("M45").Value = -2552.5908
$ws.Range("H61").Value = 11402
$ws.Range("I61").Value = 17094.727
$ws.Range("K61").Value = 17094.727
$ws.Range("M61").Value = -16882.727
$ws.Range("H97").Value = 4141.7646
$ws.Range("I97").Value = 2081.6487
$ws.Range("K97").Value = 2081.6487
$ws.Range("M97").Value = -1585.6487
$ws.Range("H110").Value = 2401.5
$ws.Range("J110").Value = 2055.5
$ws.Range("L110").Value = 2055.5
$ws.Range("N110").Value = -6145.5
$ws.Range("H116").Value = 5783.3076
$ws.Range("J116").Value = 4000
$ws.Range("L116").Value = 4000
$ws.Range("N116").Value = -8588
$ws.Range("H132").Value = 4465.7188
$ws.Range("I132").Value = 2806.1052
$ws.Range("J132").Value = 6891.3076
$ws.Range("K132").Value = 8418.3156
$ws.Range("L132").Value = 20673.9228
$ws.Range("M132").Value = -5888.3156
$ws.Range("N132").Value = -25733.9228
$ws.Range("H136").Value = 11402
$ws.Range("I136").Value = 17094.727
$ws.Range("K136").Value = 51284.181
$ws.Range("M136").Value = -48734.181

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5783.3076
$ws.Range("J3").Value = 4000
$ws.Range("L3").Value = 4000
$ws.Range("N3").Value = -4228
$ws.Range("H86").Value = 12937.5
$ws.Range("I86").Value = 15968.75
$ws.Range("J86").Value = 6875
$ws.Range("K86").Value = 15968.75
$ws.Range("L86").Value = 6875
$ws.Range("M86").Value = -14845.75
$ws.Range("N86").Value = -9121
$ws.Range("H89").Value = 12937.5
$ws.Range("I89").Value = 15968.75
$ws.Range("J89").Value = 6875
$ws.Range("K89").Value = 79843.75
$ws.Range("L89").Value = 34375
$ws.Range("M89").Value = -74227.75
$ws.Range("N89").Value = -45607

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 202034000
$ws.Range("I4").Value = 59999
$ws.Range("K4").Value = 59999
$ws.Range("M4").Value = -59887
$ws.Range("H7").Value = 431.0909
$ws.Range("I7").Value = 433.66666
$ws.Range("K7").Value = 433.66666
$ws.Range("M7").Value = -320.66666
$ws.Range("H21").Value = 30498.666
$ws.Range("I21").Value = 29999.5
$ws.Range("K21").Value = 29999.5
$ws.Range("M21").Value = -29764.5
$ws.Range("H31").Value = 3877.4167
$ws.Range("I31").Value = 3376
$ws.Range("J31").Value = 4020.6785
$ws.Range("K31").Value = 3376
$ws.Range("L31").Value = 4020.6785
$ws.Range("M31").Value = -3081
$ws.Range("N31").Value = -4610.6785
$ws.Range("H32").Value = 3033.2856
$ws.Range("I32").Value = 2705.5
$ws.Range("K32").Value = 2705.5
$ws.Range("M32").Value = -2389.5
$ws.Range("H34").Value = 3877.4167
$ws.Range("I34").Value = 3376
$ws.Range("J34").Value = 4020.6785
$ws.Range("K34").Value = 3376
$ws.Range("L34").Value = 4020.6785
$ws.Range("M34").Value = -3174
$ws.Range("N34").Value = -4424.6785
$ws.Range("H92").Value = 66549.164
$ws.Range("J92").Value = 66549.164
$ws.Range("L92").Value = 66549.164
$ws.Range("N92").Value = -71541.164
$ws.Range("H99").Value = 5933.1924
$ws.Range("J99").Value = 7280.2
$ws.Range("L99").Value = 7280.2
$ws.Range("N99").Value = -10276.2
$ws.Range("H126").Value = 5933.1924
$ws.Range("J126").Value = 7280.2
$ws.Range("L126").Value = 21840.6
$ws.Range("N126").Value = -26780.6

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1081.7142
$ws.Range("I3").Value = 1081.7142
$ws.Range("K3").Value = 3245.1426
$ws.Range("M3").Value = -3133.1426
$ws.Range("H4").Value = 38960116
$ws.Range("I4").Value = 46470356
$ws.Range("K4").Value = 139411068
$ws.Range("M4").Value = -139410956
$ws.Range("H8").Value = 2971
$ws.Range("I8").Value = 2971
$ws.Range("K8").Value = 8913
$ws.Range("M8").Value = -8774
$ws.Range("H68").Value = 2279.6
$ws.Range("I68").Value = 1780.2
$ws.Range("J68").Value = 2779
$ws.Range("K68").Value = 5340.6
$ws.Range("L68").Value = 8337
$ws.Range("M68").Value = -4529.6
$ws.Range("N68").Value = -9959
$ws.Range("H71").Value = 2279.6
$ws.Range("I71").Value = 1780.2
$ws.Range("J71").Value = 2779
$ws.Range("K71").Value = 16021.8
$ws.Range("L71").Value = 25011
$ws.Range("M71").Value = -11965.8
$ws.Range("N71").Value = -33123
$ws.Range("H75").Value = 1379.625
$ws.Range("I75").Value = 1561.4
$ws.Range("J75").Value = 1076.6666
$ws.Range("K75").Value = 4684.200000000001
$ws.Range("L75").Value = 3229.9998
$ws.Range("M75").Value = -3686.200000000001
$ws.Range("N75").Value = -5225.9998
$ws.Range("H78").Value = 1379.625
$ws.Range("I78").Value = 1561.4
$ws.Range("J78").Value = 1076.6666
$ws.Range("K78").Value = 14052.6
$ws.Range("L78").Value = 9689.9994
$ws.Range("M78").Value = -9060.6
$ws.Range("N78").Value = -19673.9994
$ws.Range("H114").Value = 1769.6666
$ws.Range("I114").Value = 1769.6666
$ws.Range("K114").Value = 5308.9998
$ws.Range("M114").Value = -2054.9998
$ws.Range("H115").Value = 1405.2
$ws.Range("I115").Value = 1256.5
$ws.Range("K115").Value = 3769.5
$ws.Range("M115").Value = -2594.5
$ws.Range("H132").Value = 1391.4445
$ws.Range("J132").Value = 1453.4286
$ws.Range("L132").Value = 13080.8574
$ws.Range("N132").Value = -18140.8574

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = $null
$ws.Range("H70").Value = 9405
$ws.Range("I70").Value = 4879
$ws.Range("J70").Value = 19965.666
$ws.Range("K70").Value = 4879
$ws.Range("L70").Value = 19965.666
$ws.Range("M70").Value = -4609
$ws.Range("N70").Value = -20505.666
$ws.Range("H73").Value = 9405
$ws.Range("I73").Value = 4879
$ws.Range("J73").Value = 19965.666
$ws.Range("K73").Value = 4879
$ws.Range("L73").Value = 19965.666
$ws.Range("M73").Value = -3943
$ws.Range("N73").Value = -21837.666
$ws.Range("H102").Value = 46250.39
$ws.Range("J102").Value = 145056.72
$ws.Range("L102").Value = 145056.72
$ws.Range("N102").Value = -148300.72
$ws.Range("H113").Value = 1141.2858
$ws.Range("I113").Value = 1141.2858
$ws.Range("K113").Value = 1141.2858
$ws.Range("M113").Value = 1028.7142
$ws.Range("H126").Value = 9013.857
$ws.Range("J126").Value = 9791
$ws.Range("L126").Value = 29373
$ws.Range("N126").Value = -34313
$ws.Range("H132").Value = 8124.353
$ws.Range("I132").Value = 9756.615
$ws.Range("J132").Value = 2819.5
$ws.Range("K132").Value = 29269.845
$ws.Range("L132").Value = 8458.5
$ws.Range("M132").Value = -26739.845
$ws.Range("N132").Value = -13518.5

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4258.087
$ws.Range("J16").Value = 5490.1113
$ws.Range("L16").Value = 5490.1113
$ws.Range("N16").Value = -5830.1113
$ws.Range("H22").Value = 2840
$ws.Range("I22").Value = 1622
$ws.Range("K22").Value = 1622
$ws.Range("M22").Value = -1327
$ws.Range("H27").Value = 2840
$ws.Range("I27").Value = 1622
$ws.Range("K27").Value = 1622
$ws.Range("M27").Value = -1515
$ws.Range("H43").Value = 30012
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = $null
$ws.Range("H46").Value = 13882.462
$ws.Range("I46").Value = 8863.5
$ws.Range("K46").Value = 8863.5
$ws.Range("M46").Value = -8675.5
$ws.Range("H61").Value = 9629.625
$ws.Range("I61").Value = 6006.8
$ws.Range("J61").Value = 15667.667
$ws.Range("K61").Value = 6006.8
$ws.Range("L61").Value = 15667.667
$ws.Range("M61").Value = -5804.8
$ws.Range("N61").Value = -16071.667
$ws.Range("H68").Value = 4621.5713
$ws.Range("I68").Value = 4337.25
$ws.Range("K68").Value = 4337.25
$ws.Range("M68").Value = -3588.25
$ws.Range("H71").Value = 4621.5713
$ws.Range("I71").Value = 4337.25
$ws.Range("K71").Value = 21686.25
$ws.Range("M71").Value = -17942.25
$ws.Range("H82").Value = 14178
$ws.Range("I82").Value = 50251
$ws.Range("J82").Value = 3871.4285
$ws.Range("K82").Value = 50251
$ws.Range("L82").Value = 3871.4285
$ws.Range("M82").Value = -49890
$ws.Range("N82").Value = -4593.4285
$ws.Range("H85").Value = 14178
$ws.Range("I85").Value = 50251
$ws.Range("J85").Value = 3871.4285
$ws.Range("K85").Value = 50251
$ws.Range("L85").Value = 3871.4285
$ws.Range("M85").Value = -49003
$ws.Range("N85").Value = -6367.4285
$ws.Range("H87").Value = 65000
$ws.Range("J87").Value = 65000
$ws.Range("L87").Value = 65000
$ws.Range("N87").Value = -67246
$ws.Range("H90").Value = 65000
$ws.Range("J90").Value = 65000
$ws.Range("L90").Value = 195000
$ws.Range("N90").Value = -206232
$ws.Range("H113").Value = 9629.625
$ws.Range("I113").Value = 6006.8
$ws.Range("J113").Value = 15667.667
$ws.Range("K113").Value = 6006.8
$ws.Range("L113").Value = 15667.667
$ws.Range("M113").Value = -3836.8
$ws.Range("N113").Value = -20007.667
$ws.Range("H122").Value = 7500
$ws.Range("I122").Value = 7500
$ws.Range("K122").Value = 22500
$ws.Range("M122").Value = -20050
$ws.Range("H131").Value = 99999
$ws.Range("J131").Value = 99999
$ws.Range("L131").Value = 99999
$ws.Range("N131").Value = -110079
$ws.Range("H132").Value = 6201
$ws.Range("I132").Value = 6647.2383
$ws.Range("J132").Value = 5793.5654
$ws.Range("K132").Value = 19941.7149
$ws.Range("L132").Value = 17380.6962
$ws.Range("M132").Value = -17411.7149
$ws.Range("N132").Value = -22440.6962
$ws.Range("H136").Value = 3769.848
$ws.Range("I136").Value = 3684.279
$ws.Range("K136").Value = 11052.837
$ws.Range("M136").Value = -8502.837

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 12633
$ws.Range("I54").Value = 12633
$ws.Range("K54").Value = 12633
$ws.Range("M54").Value = -12113
$ws.Range("H113").Value = 1472.75
$ws.Range("I113").Value = 1246
$ws.Range("K113").Value = 3738
$ws.Range("M113").Value = -1568
$ws.Range("H136").Value = 4607.6875
$ws.Range("I136").Value = 2244.125
$ws.Range("J136").Value = 6971.25
$ws.Range("K136").Value = 6732.375
$ws.Range("L136").Value = 20913.75
$ws.Range("M136").Value = -4182.375
$ws.Range("N136").Value = -26013.75
